# Weekly price-sheet update: a new daily record for "Pepino ensalada"
# (Feria Lagunitas de Puerto Montt) is inserted above the existing row 364,
# pushing rows 364-420 down to 365-421 (dimension grows from A1:R420 to
# A1:R421). The new row carries the same fixed attributes as every other
# row in this sheet (market/region/category/classification) together with
# its own date, volume and price figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push existing row 364 (and everything below it) down by one row.
$ws.Rows.Item(364).Insert()

# Populate the newly-inserted row 364 with the new record.
$ws.Range("A364").Value = 4
$ws.Range("B364").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C364").Value = "Los Lagos"
$ws.Range("D364").Value = 45034
$ws.Range("E364").Value = 10
$ws.Range("F364").Value = 100112043
$ws.Range("G364").Value = "Pepino ensalada"
$ws.Range("H364").Value = "Sin especificar"
$ws.Range("I364").Value = "Primera"
$ws.Range("J364").Value = 400
$ws.Range("K364").Value = 13000
$ws.Range("L364").Value = 13000
$ws.Range("M364").Value = 13000
$ws.Range("N364").Value = "`$/caja 60 unidades"
$ws.Range("O364").Value = "Región de Arica y Parinacota"
$ws.Range("P364").Value = 217
$ws.Range("Q364").Value = 60
$ws.Range("R364").Value = "Hortaliza"
